# KHL stats refresh - publish files + archive (2025-11-24 11:04:37)
# Regenerates Matches_SOG with newly reported fixtures and refreshes the
# rolling Shots_HA / Shots_Summary / Meta_ext snapshot tables to the
# 2025-11-23 17:00 UTC pull.

$wb = $excel.ActiveWorkbook

# --- Matches_SOG: append newly reported fixtures ---
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

# row 312: Нефтехимик vs Салават Юлаев
$wsMatches.Cells.Item(312, 1).NumberFormat = "@"
$wsMatches.Cells.Item(312, 1).Value = "897806"
$wsMatches.Cells.Item(312, 2).Value = "2025-11-23T14:00:00"
$wsMatches.Cells.Item(312, 3).Value = "Нефтехимик"
$wsMatches.Cells.Item(312, 4).Value = "Салават Юлаев"
$wsMatches.Cells.Item(312, 5).Value = 32
$wsMatches.Cells.Item(312, 6).Value = 33
$wsMatches.Cells.Item(312, 7).Value = "khl_text"

# row 313: Барыс vs Авангард
$wsMatches.Cells.Item(313, 1).NumberFormat = "@"
$wsMatches.Cells.Item(313, 1).Value = "897805"
$wsMatches.Cells.Item(313, 2).Value = "2025-11-23T15:00:00"
$wsMatches.Cells.Item(313, 3).Value = "Барыс"
$wsMatches.Cells.Item(313, 4).Value = "Авангард"
$wsMatches.Cells.Item(313, 5).Value = 26
$wsMatches.Cells.Item(313, 6).Value = 28
$wsMatches.Cells.Item(313, 7).Value = "khl_text"

# row 314: Драконы vs Ак Барс
$wsMatches.Cells.Item(314, 1).NumberFormat = "@"
$wsMatches.Cells.Item(314, 1).Value = "897807"
$wsMatches.Cells.Item(314, 2).Value = "2025-11-23T17:00:00"
$wsMatches.Cells.Item(314, 3).Value = "Драконы"
$wsMatches.Cells.Item(314, 4).Value = "Ак Барс"
$wsMatches.Cells.Item(314, 5).Value = 23
$wsMatches.Cells.Item(314, 6).Value = 38
$wsMatches.Cells.Item(314, 7).Value = "khl_text"

# row 315: Северсталь vs Лада
$wsMatches.Cells.Item(315, 1).NumberFormat = "@"
$wsMatches.Cells.Item(315, 1).Value = "897808"
$wsMatches.Cells.Item(315, 2).Value = "2025-11-23T17:00:00"
$wsMatches.Cells.Item(315, 3).Value = "Северсталь"
$wsMatches.Cells.Item(315, 4).Value = "Лада"
$wsMatches.Cells.Item(315, 5).Value = 29
$wsMatches.Cells.Item(315, 6).Value = 23
$wsMatches.Cells.Item(315, 7).Value = "khl_text"

# --- Shots_HA: refresh as_of_utc pull + updated stats ---
$wsShotsHA = $wb.Worksheets.Item("Shots_HA")

$wsShotsHA.Cells.Item(2, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(2, 6).Value = 13
$wsShotsHA.Cells.Item(2, 11).Value = 465
$wsShotsHA.Cells.Item(2, 12).Value = 406
$wsShotsHA.Cells.Item(2, 13).Value = 35.8
$wsShotsHA.Cells.Item(2, 14).Value = 31.2
$wsShotsHA.Cells.Item(3, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(4, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(5, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(5, 6).Value = 14
$wsShotsHA.Cells.Item(5, 11).Value = 478
$wsShotsHA.Cells.Item(5, 12).Value = 410
$wsShotsHA.Cells.Item(5, 13).Value = 34.1
$wsShotsHA.Cells.Item(5, 14).Value = 29.3
$wsShotsHA.Cells.Item(6, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(7, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(7, 5).Value = 19
$wsShotsHA.Cells.Item(7, 7).Value = 608
$wsShotsHA.Cells.Item(7, 8).Value = 595
$wsShotsHA.Cells.Item(7, 9).Value = 32
$wsShotsHA.Cells.Item(7, 10).Value = 31.3
$wsShotsHA.Cells.Item(8, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(9, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(10, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(10, 5).Value = 12
$wsShotsHA.Cells.Item(10, 7).Value = 342
$wsShotsHA.Cells.Item(10, 8).Value = 416
$wsShotsHA.Cells.Item(10, 9).Value = 28.5
$wsShotsHA.Cells.Item(10, 10).Value = 34.7
$wsShotsHA.Cells.Item(11, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(11, 6).Value = 15
$wsShotsHA.Cells.Item(11, 11).Value = 360
$wsShotsHA.Cells.Item(11, 12).Value = 575
$wsShotsHA.Cells.Item(11, 13).Value = 24
$wsShotsHA.Cells.Item(11, 14).Value = 38.3
$wsShotsHA.Cells.Item(12, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(13, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(14, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(14, 5).Value = 17
$wsShotsHA.Cells.Item(14, 7).Value = 522
$wsShotsHA.Cells.Item(14, 8).Value = 581
$wsShotsHA.Cells.Item(14, 9).Value = 30.7
$wsShotsHA.Cells.Item(15, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(16, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(16, 6).Value = 18
$wsShotsHA.Cells.Item(16, 11).Value = 512
$wsShotsHA.Cells.Item(16, 12).Value = 539
$wsShotsHA.Cells.Item(16, 13).Value = 28.4
$wsShotsHA.Cells.Item(16, 14).Value = 29.9
$wsShotsHA.Cells.Item(17, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(17, 5).Value = 14
$wsShotsHA.Cells.Item(17, 7).Value = 400
$wsShotsHA.Cells.Item(17, 8).Value = 319
$wsShotsHA.Cells.Item(17, 9).Value = 28.6
$wsShotsHA.Cells.Item(18, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(19, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(20, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(21, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(22, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsHA.Cells.Item(23, 4).Value = "2025-11-23T17:00:00Z"

# --- Shots_Summary: refresh as_of_utc pull + updated stats ---
$wsShotsSummary = $wb.Worksheets.Item("Shots_Summary")

$wsShotsSummary.Cells.Item(2, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(2, 5).Value = 27
$wsShotsSummary.Cells.Item(2, 6).Value = 914
$wsShotsSummary.Cells.Item(2, 7).Value = 810
$wsShotsSummary.Cells.Item(2, 8).Value = 33.9
$wsShotsSummary.Cells.Item(2, 9).Value = 30
$wsShotsSummary.Cells.Item(3, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(4, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(5, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(5, 5).Value = 30
$wsShotsSummary.Cells.Item(5, 6).Value = 1006
$wsShotsSummary.Cells.Item(5, 7).Value = 818
$wsShotsSummary.Cells.Item(5, 8).Value = 33.5
$wsShotsSummary.Cells.Item(5, 9).Value = 27.3
$wsShotsSummary.Cells.Item(6, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(7, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(7, 5).Value = 30
$wsShotsSummary.Cells.Item(7, 6).Value = 923
$wsShotsSummary.Cells.Item(7, 7).Value = 962
$wsShotsSummary.Cells.Item(7, 8).Value = 30.8
$wsShotsSummary.Cells.Item(7, 9).Value = 32.1
$wsShotsSummary.Cells.Item(8, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(9, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(10, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(10, 5).Value = 28
$wsShotsSummary.Cells.Item(10, 6).Value = 760
$wsShotsSummary.Cells.Item(10, 7).Value = 986
$wsShotsSummary.Cells.Item(10, 8).Value = 27.1
$wsShotsSummary.Cells.Item(10, 9).Value = 35.2
$wsShotsSummary.Cells.Item(11, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(11, 5).Value = 28
$wsShotsSummary.Cells.Item(11, 6).Value = 713
$wsShotsSummary.Cells.Item(11, 7).Value = 1042
$wsShotsSummary.Cells.Item(11, 8).Value = 25.5
$wsShotsSummary.Cells.Item(11, 9).Value = 37.2
$wsShotsSummary.Cells.Item(12, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(13, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(14, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(14, 5).Value = 30
$wsShotsSummary.Cells.Item(14, 6).Value = 896
$wsShotsSummary.Cells.Item(14, 7).Value = 1056
$wsShotsSummary.Cells.Item(14, 8).Value = 29.9
$wsShotsSummary.Cells.Item(14, 9).Value = 35.2
$wsShotsSummary.Cells.Item(15, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(16, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(16, 5).Value = 29
$wsShotsSummary.Cells.Item(16, 6).Value = 804
$wsShotsSummary.Cells.Item(16, 7).Value = 853
$wsShotsSummary.Cells.Item(16, 8).Value = 27.7
$wsShotsSummary.Cells.Item(16, 9).Value = 29.4
$wsShotsSummary.Cells.Item(17, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(17, 5).Value = 29
$wsShotsSummary.Cells.Item(17, 6).Value = 893
$wsShotsSummary.Cells.Item(17, 7).Value = 721
$wsShotsSummary.Cells.Item(17, 8).Value = 30.8
$wsShotsSummary.Cells.Item(18, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(19, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(20, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(21, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(22, 4).Value = "2025-11-23T17:00:00Z"
$wsShotsSummary.Cells.Item(23, 4).Value = "2025-11-23T17:00:00Z"

# --- Meta_ext: refresh as_of_utc pull + updated stats ---
$wsMeta = $wb.Worksheets.Item("Meta_ext")

$wsMeta.Cells.Item(2, 2).Value = "2025-11-23T17:00:00Z"
$wsMeta.Cells.Item(2, 4).Value = 9

